$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cell values (rows 2-6) ---
$ws.Range("H2").Value = 3.35
$ws.Range("I2").Value = 3.5
$ws.Range("P2").Value = 2.74
$ws.Range("Q2").Value = 1.52
$ws.Range("S2").Value = 2.26
$ws.Range("T2").Value = 1.51
$ws.Range("U2").Value = 2.68
$ws.Range("Y2").Value = 22.0
$ws.Range("AA2").Value = 980.0
$ws.Range("AB2").Value = 16.5
$ws.Range("AC2").Value = 11.0
$ws.Range("AD2").Value = 16.0
$ws.Range("AG2").Value = 12.0
$ws.Range("AH2").Value = 15.5
$ws.Range("AI2").Value = 36.0
$ws.Range("AJ2").Value = 1000.0
$ws.Range("AN2").Value = 9.8
$ws.Range("F3").Value = 2.0
$ws.Range("G3").Value = 2.16
$ws.Range("H3").Value = 3.45
$ws.Range("I3").Value = 3.9
$ws.Range("J3").Value = 3.8
$ws.Range("K3").Value = 4.3
$ws.Range("O3").Value = 1.22
$ws.Range("P3").Value = 2.34
$ws.Range("Q3").Value = 1.64
$ws.Range("R3").Value = 1.53
$ws.Range("S3").Value = 2.6
$ws.Range("T3").Value = 1.61
$ws.Range("U3").Value = 2.4
$ws.Range("V3").Value = 1.34
$ws.Range("W3").Value = 1.87
$ws.Range("X3").Value = 26.0
$ws.Range("F4").Value = 5.5
$ws.Range("I4").Value = 1.75
$ws.Range("K4").Value = 4.4
$ws.Range("N4").Value = 3.8
$ws.Range("P4").Value = 1.97
$ws.Range("Q4").Value = 1.84
$ws.Range("R4").Value = 1.37
$ws.Range("T4").Value = 1.83
$ws.Range("U4").Value = 2.0
$ws.Range("V4").Value = 2.32
$ws.Range("X4").Value = 19.5
$ws.Range("Y4").Value = 10.5
$ws.Range("AA4").Value = 21.0
$ws.Range("AC4").Value = 11.5
$ws.Range("AE4").Value = 22.0
$ws.Range("AF4").Value = 60.0
$ws.Range("AG4").Value = 26.0
$ws.Range("AI4").Value = 42.0
$ws.Range("AK4").Value = 100.0
$ws.Range("AO4").Value = 12.0
$ws.Range("F6").Value = 2.96
$ws.Range("N6").Value = 2.58
$ws.Range("O6").Value = 1.55
$ws.Range("V6").Value = 1.53
$ws.Range("AF6").Value = 980.0

# --- Add new rows 7, 8, 9 ---

# Row 7
$ws.Range("A7").Value = 'Colombian Primera B'
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = '2025-10-16'
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Value = '21:00:00'
$ws.Range("D7").Value = 'Quindio'
$ws.Range("E7").Value = 'Tigres FC Zipaquira'
$ws.Range("F7").Value = 1.57
$ws.Range("G7").Value = 1.62
$ws.Range("H7").Value = 6.4
$ws.Range("I7").Value = 9.8
$ws.Range("J7").Value = 3.7
$ws.Range("K7").Value = 4.4
$ws.Range("L7").Value = 1.01
$ws.Range("M7").Value = 1.08
$ws.Range("N7").Value = 1.71
$ws.Range("O7").Value = 1.39
$ws.Range("P7").Value = 1.71
$ws.Range("Q7").Value = 2.12
$ws.Range("R7").Value = 1.26
$ws.Range("S7").Value = 3.95
$ws.Range("T7").Value = 2.16
$ws.Range("U7").Value = 1.7
$ws.Range("V7").Value = 1.12
$ws.Range("W7").Value = 2.6
$ws.Range("X7").Value = 12.5
$ws.Range("Y7").Value = 21.0
$ws.Range("Z7").Value = 70.0
$ws.Range("AA7").Value = 310.0
$ws.Range("AB7").Value = 6.8
$ws.Range("AC7").Value = 9.4
$ws.Range("AD7").Value = 32.0
$ws.Range("AE7").Value = 170.0
$ws.Range("AF7").Value = 8.8
$ws.Range("AG7").Value = 11.0
$ws.Range("AH7").Value = 38.0
$ws.Range("AI7").Value = 240.0
$ws.Range("AJ7").Value = 15.0
$ws.Range("AK7").Value = 20.0
$ws.Range("AL7").Value = 55.0
$ws.Range("AM7").Value = 230.0
$ws.Range("AN7").Value = 13.0
$ws.Range("AO7").Value = 320.0

# Row 8
$ws.Range("A8").Value = 'Brazilian Serie A'
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = '2025-10-16'
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = '21:30:00'
$ws.Range("D8").Value = 'EC Vitoria Salvador'
$ws.Range("E8").Value = 'Bahia'
$ws.Range("F8").Value = 2.9
$ws.Range("G8").Value = 3.15
$ws.Range("H8").Value = 2.72
$ws.Range("I8").Value = 2.94
$ws.Range("J8").Value = 3.15
$ws.Range("K8").Value = 3.2
$ws.Range("L8").Value = 1.57
$ws.Range("M8").Value = 1.12
$ws.Range("N8").Value = 2.62
$ws.Range("O8").Value = 1.54
$ws.Range("P8").Value = 1.52
$ws.Range("Q8").Value = 2.62
$ws.Range("R8").Value = 1.19
$ws.Range("S8").Value = 5.4
$ws.Range("T8").Value = 2.08
$ws.Range("U8").Value = 1.77
$ws.Range("V8").Value = 1.52
$ws.Range("W8").Value = 1.47
$ws.Range("X8").Value = 9.2
$ws.Range("Y8").Value = 8.4
$ws.Range("Z8").Value = 17.0
$ws.Range("AA8").Value = 55.0
$ws.Range("AB8").Value = 8.8
$ws.Range("AC8").Value = 7.0
$ws.Range("AD8").Value = 14.0
$ws.Range("AE8").Value = 42.0
$ws.Range("AF8").Value = 18.0
$ws.Range("AG8").Value = 15.0
$ws.Range("AH8").Value = 23.0
$ws.Range("AI8").Value = 70.0
$ws.Range("AJ8").Value = 60.0
$ws.Range("AK8").Value = 48.0
$ws.Range("AL8").Value = 95.0
$ws.Range("AM8").Value = 170.0
$ws.Range("AN8").Value = 55.0
$ws.Range("AO8").Value = 50.0

# Row 9
$ws.Range("A9").Value = 'Brazilian Serie A'
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = '2025-10-16'
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = '21:30:00'
$ws.Range("D9").Value = 'Fluminense'
$ws.Range("E9").Value = 'Juventude'
$ws.Range("F9").Value = 1.37
$ws.Range("G9").Value = 1.41
$ws.Range("H9").Value = 11.5
$ws.Range("I9").Value = 13.5
$ws.Range("J9").Value = 4.8
$ws.Range("K9").Value = 5.4
$ws.Range("L9").Value = 1.44
$ws.Range("M9").Value = 1.07
$ws.Range("N9").Value = 3.3
$ws.Range("O9").Value = 1.37
$ws.Range("P9").Value = 1.8
$ws.Range("Q9").Value = 2.1
$ws.Range("R9").Value = 1.29
$ws.Range("S9").Value = 3.8
$ws.Range("T9").Value = 2.5
$ws.Range("U9").Value = 1.56
$ws.Range("V9").Value = 1.08
$ws.Range("W9").Value = 3.35
$ws.Range("X9").Value = 14.5
$ws.Range("Y9").Value = 30.0
$ws.Range("Z9").Value = 140.0
$ws.Range("AA9").Value = 1000.0
$ws.Range("AB9").Value = 6.6
$ws.Range("AC9").Value = 12.5
$ws.Range("AD9").Value = 980.0
$ws.Range("AE9").Value = 350.0
$ws.Range("AF9").Value = 7.2
$ws.Range("AG9").Value = 11.5
$ws.Range("AH9").Value = 980.0
$ws.Range("AI9").Value = 280.0
$ws.Range("AJ9").Value = 11.0
$ws.Range("AK9").Value = 19.5
$ws.Range("AL9").Value = 65.0
$ws.Range("AM9").Value = 400.0
$ws.Range("AN9").Value = 8.6
$ws.Range("AO9").Value = 1000.0
